$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1188.1333
$ws.Range("I19").Value = 1351.6
$ws.Range("K19").Value = 1351.6
$ws.Range("M19").Value = -1176.6

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 712.25
$ws.Range("I39").Value = 116.5
$ws.Range("K39").Value = 349.5
$ws.Range("M39").Value = -53.5

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 39.125
$ws.Range("I41").Value = 45
$ws.Range("K41").Value = 45
$ws.Range("M41").Value = 395

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 10014.96
$ws.Range("J69").Value = 10014.96
$ws.Range("L69").Value = 30044.88
$ws.Range("N69").Value = -31792.88

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 10014.96
$ws.Range("J72").Value = 10014.96
$ws.Range("L72").Value = 90134.63999999998
$ws.Range("N72").Value = -98870.63999999998

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 8181.6665
$ws.Range("I125").Value = 600
$ws.Range("K125").Value = 5400
$ws.Range("M125").Value = -2940

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3564.1667
$ws.Range("I131").Value = 2496.3333
$ws.Range("J131").Value = 5699.8335
$ws.Range("K131").Value = 7488.999899999999
$ws.Range("L131").Value = 17099.5005
$ws.Range("M131").Value = -2448.999899999999
$ws.Range("N131").Value = -27179.5005

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3666.1738
$ws.Range("I61").Value = 3316.1177
$ws.Range("K61").Value = 3316.1177
$ws.Range("M61").Value = -3104.1177

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5874.1353
$ws.Range("I122").Value = 4769.1113
$ws.Range("K122").Value = 14307.3339
$ws.Range("M122").Value = -11857.3339

# ARM row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 115499.5
$ws.Range("J131").Value = 115499.5
$ws.Range("L131").Value = 115499.5
$ws.Range("N131").Value = -125579.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3666.1738
$ws.Range("I136").Value = 3316.1177
$ws.Range("K136").Value = 9948.3531
$ws.Range("M136").Value = -7398.3531

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1069.6
$ws.Range("I86").Value = 1162
$ws.Range("J86").Value = 700
$ws.Range("K86").Value = 1162
$ws.Range("L86").Value = 700
$ws.Range("M86").Value = -39
$ws.Range("N86").Value = -2946

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1069.6
$ws.Range("I89").Value = 1162
$ws.Range("J89").Value = 700
$ws.Range("K89").Value = 5810
$ws.Range("L89").Value = 3500
$ws.Range("M89").Value = -194
$ws.Range("N89").Value = -14732

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1123.3334
$ws.Range("I107").Value = 763.75
$ws.Range("K107").Value = 763.75
$ws.Range("M107").Value = 1156.25

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1485.8695
$ws.Range("I16").Value = 1536.875
$ws.Range("K16").Value = 1536.875
$ws.Range("M16").Value = -1249.875

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1485.8695
$ws.Range("I113").Value = 1536.875
$ws.Range("K113").Value = 1536.875
$ws.Range("M113").Value = 633.125

# CRP row 118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 105330
$ws.Range("J118").Value = 105330
$ws.Range("L118").Value = 105330
$ws.Range("N118").Value = -108644

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1916.6
$ws.Range("I18").Value = 395.2857
$ws.Range("J18").Value = 5466.3335
$ws.Range("K18").Value = 1185.8571
$ws.Range("L18").Value = 16399.0005
$ws.Range("M18").Value = -1016.8571
$ws.Range("N18").Value = -16737.0005

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1530.375
$ws.Range("J131").Value = 1893.9231
$ws.Range("L131").Value = 5681.7693
$ws.Range("N131").Value = -15761.7693

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1872.3043
$ws.Range("I102").Value = 1766.5264
$ws.Range("J102").Value = 2374.75
$ws.Range("K102").Value = 1766.5264
$ws.Range("L102").Value = 2374.75
$ws.Range("M102").Value = -144.5264
$ws.Range("N102").Value = -5618.75

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3210.2144
$ws.Range("I126").Value = 2999.3
$ws.Range("J126").Value = 3737.5
$ws.Range("K126").Value = 8997.900000000001
$ws.Range("L126").Value = 11212.5
$ws.Range("M126").Value = -6527.900000000001
$ws.Range("N126").Value = -16152.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2549.2156
$ws.Range("J132").Value = 3656.077
$ws.Range("L132").Value = 10968.231
$ws.Range("N132").Value = -16028.231

# LTW row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 18643.5
$ws.Range("I24").Value = 14858.333
$ws.Range("K24").Value = 14858.333
$ws.Range("M24").Value = -14515.333

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2159.8
$ws.Range("I40").Value = 2099.25
$ws.Range("J40").Value = 2402
$ws.Range("K40").Value = 2099.25
$ws.Range("L40").Value = 2402
$ws.Range("M40").Value = -1963.25
$ws.Range("N40").Value = -2674

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4499.3335
$ws.Range("I61").Value = 3399.2
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 3399.2
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -3197.2
$ws.Range("N61").Value = -10404

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4499.3335
$ws.Range("I113").Value = 3399.2
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 3399.2
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -1229.2
$ws.Range("N113").Value = -14340

# LTW row 116
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 107116.336
$ws.Range("J116").Value = 107116.336
$ws.Range("L116").Value = 107116.336
$ws.Range("N116").Value = -116294.336

# WVR row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 28678
$ws.Range("I21").Value = 28000
$ws.Range("J21").Value = 28813.6
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 28813.6
$ws.Range("M21").Value = -27765
$ws.Range("N21").Value = -29283.6

# WVR row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 28678
$ws.Range("I35").Value = 28000
$ws.Range("J35").Value = 28813.6
$ws.Range("K35").Value = 28000
$ws.Range("L35").Value = 28813.6
$ws.Range("M35").Value = -27710
$ws.Range("N35").Value = -29393.6

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 81526.62
$ws.Range("I81").Value = 253214.25
$ws.Range("J81").Value = 5221
$ws.Range("K81").Value = 506428.5
$ws.Range("L81").Value = 10442
$ws.Range("M81").Value = -505367.5
$ws.Range("N81").Value = -12564

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 81526.62
$ws.Range("I84").Value = 253214.25
$ws.Range("J84").Value = 5221
$ws.Range("K84").Value = 2532142.5
$ws.Range("L84").Value = 52210
$ws.Range("M84").Value = -2526838.5
$ws.Range("N84").Value = -62818
